$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Update the pincode value in A4 (560001 -> 111)
$ws.Range("A4").Value = 111

# Add vertical centering to cells that only had horizontal centering so they
# merge with the existing "center/center" styles already used by A1 and A2.
$ws.Range("B1:C1").VerticalAlignment = -4108
$ws.Range("A2:C4").VerticalAlignment = -4108

# Update the saved selection to C7
[void]$ws.Range("C7").Select()
